# Apio (Vega Modelo de Temuco) weekly update:
# a new price-report row is inserted at row 125, pushing every
# subsequent row (old 125..243) down by one (new 126..244).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 125 (shifts 125..243 -> 126..244)
$ws.Rows.Item(125).EntireRow.Insert()

# Populate the newly inserted row 125 with the new weekly record
$ws.Cells.Item(125, 1).Value  = 10
$ws.Cells.Item(125, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(125, 3).Value  = "La Araucanía"
$ws.Cells.Item(125, 4).Value  = 44589
$ws.Cells.Item(125, 5).Value  = 9
$ws.Cells.Item(125, 6).Value  = 100112017
$ws.Cells.Item(125, 7).Value  = "Apio"
$ws.Cells.Item(125, 8).Value  = "Americana (o)"
$ws.Cells.Item(125, 9).Value  = "Primera"
$ws.Cells.Item(125, 10).Value = 50
$ws.Cells.Item(125, 11).Value = 10000
$ws.Cells.Item(125, 12).Value = 10000
$ws.Cells.Item(125, 13).Value = 10000
$ws.Cells.Item(125, 14).Value = "`$/docena de matas"
$ws.Cells.Item(125, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(125, 16).Value = 1667
$ws.Cells.Item(125, 17).Value = 6
$ws.Cells.Item(125, 18).Value = "Hortaliza"
